# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Agrícola del Norte S.A. de Arica - Mango"
# at row 65, pushing the existing rows 65-121 down to 66-122, and populate the
# newly inserted row with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 65..121 down by one row (Excel copies formatting from the row above).
$ws.Rows.Item(65).Insert()

# Populate the newly-inserted row 65 with the new week's record.
$ws.Cells.Item(65, 1).Value = 1
$ws.Cells.Item(65, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(65, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(65, 4).Value = 44586
$ws.Cells.Item(65, 5).Value = 15
$ws.Cells.Item(65, 6).Value = 'Fruta'
$ws.Cells.Item(65, 7).Value = 100108
$ws.Cells.Item(65, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(65, 9).Value = 100108002
$ws.Cells.Item(65, 10).Value = 'Mango'
$ws.Cells.Item(65, 11).Value = 'Sin especificar'
$ws.Cells.Item(65, 12).Value = 'Especial'
$ws.Cells.Item(65, 13).Value = 450
$ws.Cells.Item(65, 14).Value = 5000
$ws.Cells.Item(65, 15).Value = 5500
$ws.Cells.Item(65, 16).Value = 5250
$ws.Cells.Item(65, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(65, 18).Value = 'Perú'
$ws.Cells.Item(65, 19).Value = 1312
$ws.Cells.Item(65, 20).Value = 4
